# Edit: add weekly attendance columns G (9.16) and H (9.22) to the
# "24网络技术3班考勤" attendance sheet, matching the source commit's XLSX diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$check = [char]0x221A   # "√" (present)
$cross = [char]0x00D7   # "×" (absent)

# Header row: two new date columns appended after F (9.15).
$ws.Cells.Item(2, 7).Value = 9.16
$ws.Cells.Item(2, 8).Value = 9.22

# Per-student attendance marks for the two new columns (row, G-mark, H-mark).
# "C" = check (present), "X" = cross (absent).
$attendance = @(
    @(3,"C","C"),
    @(4,"C","C"),
    @(5,"X","X"),
    @(6,"C","C"),
    @(7,"C","C"),
    @(8,"C","C"),
    @(9,"C","C"),
    @(10,"C","C"),
    @(11,"X","X"),
    @(12,"C","C"),
    @(13,"C","C"),
    @(14,"C","C"),
    @(15,"C","C"),
    @(16,"X","X"),
    @(17,"X","C"),
    @(18,"C","C"),
    @(19,"C","X"),
    @(20,"C","C"),
    @(21,"C","C"),
    @(22,"C","C"),
    @(23,"C","C"),
    @(24,"X","C"),
    @(25,"C","X"),
    @(26,"C","X"),
    @(27,"X","X"),
    @(28,"C","C"),
    @(29,"C","C"),
    @(30,"C","C"),
    @(31,"X","X"),
    @(32,"X","C"),
    @(33,"C","C"),
    @(34,"X","C"),
    @(35,"C","C"),
    @(36,"C","C"),
    @(37,"C","C"),
    @(38,"X","X"),
    @(39,"C","C"),
    @(40,"C","C"),
    @(41,"C","C"),
    @(42,"C","C"),
    @(43,"C","C"),
    @(44,"C","C"),
    @(45,"X","X"),
    @(46,"C","C"),
    @(47,"C","X"),
    @(48,"C","C"),
    @(49,"C","C")
)

foreach ($entry in $attendance) {
    $r = $entry[0]
    if ($entry[1] -eq "C") { $gSym = $check } else { $gSym = $cross }
    if ($entry[2] -eq "C") { $hSym = $check } else { $hSym = $cross }
    $ws.Cells.Item($r, 7).Value = $gSym
    $ws.Cells.Item($r, 8).Value = $hSym
}

# Match the saved selection/scroll state from the edit.
$ws.Range("H12").Select()

Write-Host "Attendance columns G and H populated for rows 3-49."
